# team1_aufgabe4_lernjournal.xlsx - "update lernjournal 4 - change times"
#
# The three daily "Tätigkeiten" blocks (rows 11-13, 23-25, 35-37) get their
# three entries reshuffled (a cyclic rotation) and the dates on the two
# "Bearbeitung" rows updated to match the new day.
#
# New per-block layout (row offset 0/1/2 from the block's first data row):
#   offset 0: "Besprechung Konzepte Aufgabe 4" | 50  min | 21.11.2020
#   offset 1: "Bearbeitung Aufgabe 4"          | 150 min | 22.11.2020
#   offset 2: "Bearbeitung Präsentation"       | 30  min | 22.11.2020

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Serial date numbers (Excel 1900 date system)
$dNov21 = 44156   # 21.11.2020
$dNov22 = 44157   # 22.11.2020

$blocks = @(11, 23, 35)

foreach ($firstRow in $blocks) {
    $r0 = $firstRow
    $r1 = $firstRow + 1
    $r2 = $firstRow + 2

    # Row 0: Besprechung Konzepte Aufgabe 4 - 50 min - 21.11.2020
    $ws.Cells.Item($r0, 1).Value = "Besprechung Konzepte Aufgabe 4"
    $ws.Cells.Item($r0, 2).Value = 50
    $ws.Cells.Item($r0, 4).Value = $dNov21

    # Row 1: Bearbeitung Aufgabe 4 - 150 min - 22.11.2020
    $ws.Cells.Item($r1, 1).Value = "Bearbeitung Aufgabe 4"
    $ws.Cells.Item($r1, 2).Value = 150
    $ws.Cells.Item($r1, 4).Value = $dNov22

    # Row 2: Bearbeitung Präsentation - 30 min - 22.11.2020
    $ws.Cells.Item($r2, 1).Value = "Bearbeitung Präsentation"
    $ws.Cells.Item($r2, 2).Value = 30
    $ws.Cells.Item($r2, 4).Value = $dNov22
}

# Restore the active selection to match the final saved state of the sheet
$ws.Range("F16").Select()
